{"js": "// Delete the \"Requisitos\" section (its Heading2 title paragraph and the\n// following \"LOQ4205 - ... (Requisito fraco)\" list-bullet paragraph) from\n// the end of the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = items.length - 1; i >= 0; i--) {\n  const style = (items[i].style || \"\").trim();\n  const text = (items[i].text || \"\").trim();\n  const isRequisitosHeading = style === \"Heading 2\" && text === \"Requisitos\";\n  const isRequisitoItem = text.indexOf(\"LOQ4205\") === 0 && text.indexOf(\"Requisito fraco\") !== -1;\n  if (isRequisitosHeading || isRequisitoItem) {\n    items[i].delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Remove the \"Requisitos\" section from the end of the document: the\n# Heading 2 paragraph titled \"Requisitos\" and the following List Bullet\n# paragraph \"LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)\".\n$d = $word.ActiveDocument\n\n$toDelete = @()\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    $text = $p.Range.Text.Trim()\n\n    $isRequisitosHeading = ($styleName -eq \"Heading 2\") -and ($text -eq \"Requisitos\")\n    $isRequisitoItem = $text.StartsWith(\"LOQ4205\") -and $text.Contains(\"Requisito fraco\")\n\n    if ($isRequisitosHeading -or $isRequisitoItem) {\n        $toDelete += $i\n    }\n}\n\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
